# Расходники 9 октября.xlsx - bot update:
#  - auto-resize / explicit column widths on Sheet1
#  - G column ("ImageURL") gets an explicit 0 instead of being blank
#  - a handful of quantity corrections in column F
#  - stray empty "Модель" / "Характеристика" cells (leftover blank strings) cleared

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (auto-resize columns) -----------------------------------
# ColumnWidth is specified in characters; the engine stores width in the xlsx
# with Excel's standard +5px (0.8333... chars) padding, so subtract that back
# out to land on the exact target stored widths (4, 45, 29, 24, 44, 12, 10).
$padding = 0.8333333333333334
$ws.Columns.Item(1).ColumnWidth = 4 - $padding
$ws.Columns.Item(2).ColumnWidth = 45 - $padding
$ws.Columns.Item(3).ColumnWidth = 29 - $padding
$ws.Columns.Item(4).ColumnWidth = 24 - $padding
$ws.Columns.Item(5).ColumnWidth = 44 - $padding
$ws.Columns.Item(6).ColumnWidth = 12 - $padding
$ws.Columns.Item(7).ColumnWidth = 10 - $padding

# --- Quantity corrections in column F ---------------------------------------
$ws.Range("F5").Value = 3
$ws.Range("F9").Value = 15
$ws.Range("F11").Value = 3
$ws.Range("F21").Value = 10

# --- Column G ("ImageURL") defaults to 0 instead of an empty string ---------
for ($row = 2; $row -le 62; $row++) {
    $ws.Cells.Item($row, 7).Value = 0
}

# --- Clear stray empty cells left over in "Модель" (C) / "Характеристика" (E)
$emptyModelCells = @("C7", "C8", "C20", "C28", "C29", "C44", "C45", "C46", "C60")
foreach ($cellRef in $emptyModelCells) {
    $ws.Range($cellRef).ClearContents()
}
$ws.Range("E60").ClearContents()
